# Applies: "Thêm phân hệ nhân viên và vẽ màn phân quyền"
# Adds two new mini table-definitions to Sheet1:
#   - "Nhân viên" (Employee) table at rows 44-46
#   - "Quyền" (Rules) table at rows 48-50
# following the exact same layout pattern used by the other table
# definitions already present on the sheet (a colored title row with the
# Vietnamese table name + technical table name, followed by one row of
# Vietnamese field names and one row of the matching field codes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlFormats = $xlPasteFormats

function Copy-Format($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial($xlFormats) | Out-Null
}

# ---------------------------------------------------------------------
# "Nhân viên" (Employee) table - rows 44-46
# ---------------------------------------------------------------------

# Title row (A44:B44) - reuse the look of the other table-title rows:
# A = highlighted Vietnamese table name (plain yellow fill, no border),
# B = technical/table code name (grey "code" font, same as other tables).
Copy-Format "B21" "B44"
$ws.Range("A44").Value = "Nhân viên"
$ws.Range("A44").Interior.Color = 65535
$ws.Range("B44").Value = "ldt_employee"

# Field header rows (plain bordered cells, same style used throughout
# the sheet for the field rows - e.g. row 45/46 below a title row).
Copy-Format "A22" "A45"
Copy-Format "B22" "B45"
Copy-Format "C22" "C45"
Copy-Format "D22" "D45"
$ws.Range("A45").Value = "ID Nhân viên"
$ws.Range("B45").Value = "Tên nhân viên"
$ws.Range("C45").Value = "Mã nhân viên"
$ws.Range("D45").Value = "Quyền"

Copy-Format "A23" "A46"
Copy-Format "B23" "B46"
Copy-Format "C23" "C46"
Copy-Format "D23" "D46"
$ws.Range("A46").Value = "EmployeeID"
$ws.Range("B46").Value = "EmployeeName"
$ws.Range("C46").Value = "EmployeeCode"
$ws.Range("D46").Value = "Rules"

# ---------------------------------------------------------------------
# "Quyền" (Rules) table - rows 48-50
# ---------------------------------------------------------------------

# Title row (A48:B48)
Copy-Format "A8" "A48"
Copy-Format "B8" "B48"
$ws.Range("A48").Value = "Quyền"
$ws.Range("B48").Value = "ldt_rule"

Copy-Format "A9" "A49"
Copy-Format "B9" "B49"
$ws.Range("A49").Value = "ID quyền"
$ws.Range("B49").Value = "Tên quyền"

Copy-Format "A10" "A50"
Copy-Format "B10" "B50"
$ws.Range("A50").Value = "RuleID"
$ws.Range("B50").Value = "RuleName"

# Restore selection/active cell similar to the final state left in Excel.
$excel.CutCopyMode = $false
$ws.Range("D46").Select() | Out-Null
